$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").Value = 6
$ws.Range("E1").Value = 12

$ws.Range("D2").Value = 3.25
$ws.Range("E2").Value = 3.75

$ws.Range("D3").Value = 4.25
$ws.Range("E3").Value = 6.5

$ws.Range("C4").Value = 3.75
$ws.Range("D4").Value = 3.5

$ws.Range("C5").Value = 1.83
$ws.Range("D5").Value = 3.6
$ws.Range("E5").Value = 4.3

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 3.4
$ws.Range("E6").Value = 2.35

$ws.Range("C7").Value = 2.5499999999999998
$ws.Range("E7").Value = 2.7

$ws.Range("C8").Value = 2.15
$ws.Range("E8").Value = 3.25

$ws.Range("C9").Value = 1.33
$ws.Range("D9").Value = 5.5
$ws.Range("E9").Value = 8.5

$ws.Range("C10").Value = 2.2000000000000002
$ws.Range("E10").Value = 3.3

$ws.Range("E8").Select()
